# Page Object Model (POM) added: populate row 2 with a second login
# (email/password) pair, matching the existing admin@gmail.com /
# password layout on row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values -> become new shared-string entries (indices 1 and 2,
# alongside the existing "admin@gmail.com" at index 0).
$ws.Range("A2").Value = "siloshasilu@gmail.com"
$ws.Range("B2").Value = "November@1"

# Hyperlink both new cells the same way A1 already is (mailto: link,
# display text equal to the cell text).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:siloshasilu@gmail.com", "", "", "siloshasilu@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:November@1", "", "", "November@1")

# Hyperlinks.Add auto-applies the built-in blue/underlined "Hyperlink"
# style; reset the new cells back to the same plain blue font already
# used by A1 (no underline, no theme/hyperlink style) so A2/B2 share
# A1's cell style instead of creating a new one.
$ws.Range("A2:B2").Font.Name = "Arial"
$ws.Range("A2:B2").Font.Underline = $excel.XlUnderlineStyle.xlUnderlineStyleNone
$ws.Range("A2:B2").Font.Color = $ws.Range("A1").Font.Color

# Matches the sheet's new selection/active cell after the edit.
$ws.Range("B2").Select()
